$wb = $excel.ActiveWorkbook
$wsTrials = $wb.Worksheets.Item("Trials")
$wsSummary = $wb.Worksheets.Item("Summary")

$wsTrials.Range("B202").Value = 1537
$wsTrials.Range("C202").Value = 0.7718706130981445
$wsTrials.Range("B203").Value = 1748
$wsTrials.Range("C203").Value = 0.8732006549835205
$wsTrials.Range("B204").Value = 2906
$wsTrials.Range("C204").Value = 1.624788045883179
$wsTrials.Range("B205").Value = 9435
$wsTrials.Range("C205").Value = 9.90837287902832
$wsTrials.Range("B206").Value = 3980
$wsTrials.Range("C206").Value = 2.881757497787476
$wsTrials.Range("B207").Value = 9161
$wsTrials.Range("C207").Value = 10.11763405799866
$wsTrials.Range("B208").Value = 4384
$wsTrials.Range("C208").Value = 3.039589643478394
$wsTrials.Range("B209").Value = 1368
$wsTrials.Range("C209").Value = 0.8540494441986084
$wsTrials.Range("B210").Value = 2853
$wsTrials.Range("C210").Value = 1.882181882858276
$wsTrials.Range("B211").Value = 11676
$wsTrials.Range("C211").Value = 15.58777832984924
$wsTrials.Range("B212").Value = 3812
$wsTrials.Range("C212").Value = 2.715559005737305
$wsTrials.Range("B213").Value = 19975
$wsTrials.Range("C213").Value = 43.01579236984253
$wsTrials.Range("B214").Value = 2474
$wsTrials.Range("C214").Value = 1.634118318557739
$wsTrials.Range("B215").Value = 5324
$wsTrials.Range("C215").Value = 4.47121000289917
$wsTrials.Range("B216").Value = 5457
$wsTrials.Range("C216").Value = 4.495755434036255
$wsTrials.Range("B217").Value = 10629
$wsTrials.Range("C217").Value = 13.41584444046021
$wsTrials.Range("B218").Value = 2528
$wsTrials.Range("C218").Value = 1.444885730743408
$wsTrials.Range("B219").Value = 3127
$wsTrials.Range("C219").Value = 1.333049774169922
$wsTrials.Range("B220").Value = 3061
$wsTrials.Range("C220").Value = 1.878376483917236
$wsTrials.Range("B221").Value = 9802
$wsTrials.Range("C221").Value = 10.47914052009583
$wsTrials.Range("B222").Value = 11050
$wsTrials.Range("C222").Value = 12.16622686386108
$wsTrials.Range("B223").Value = 8363
$wsTrials.Range("C223").Value = 7.687073707580566
$wsTrials.Range("B224").Value = 1089
$wsTrials.Range("C224").Value = 0.295243501663208
$wsTrials.Range("B225").Value = 9478
$wsTrials.Range("C225").Value = 9.774614334106445
$wsTrials.Range("D225").Value = 1
$wsTrials.Range("B226").Value = 3351
$wsTrials.Range("C226").Value = 2.010892629623413
$wsTrials.Range("B227").Value = 9283
$wsTrials.Range("C227").Value = 9.715269565582275
$wsTrials.Range("B228").Value = 23555
$wsTrials.Range("C228").Value = 50.18421578407288
$wsTrials.Range("B229").Value = 4098
$wsTrials.Range("C229").Value = 2.361819505691528
$wsTrials.Range("B230").Value = 6334
$wsTrials.Range("C230").Value = 4.515548944473267
$wsTrials.Range("B231").Value = 4002
$wsTrials.Range("C231").Value = 2.453101396560669
$wsTrials.Range("B232").Value = 8028
$wsTrials.Range("C232").Value = 7.034000158309937
$wsTrials.Range("B233").Value = 511
$wsTrials.Range("C233").Value = 0.1745448112487793
$wsTrials.Range("B234").Value = 4305
$wsTrials.Range("C234").Value = 2.689745664596558
$wsTrials.Range("B235").Value = 4993
$wsTrials.Range("C235").Value = 3.431664228439331
$wsTrials.Range("B236").Value = 706
$wsTrials.Range("C236").Value = 0.2728095054626465
$wsTrials.Range("B237").Value = 8314
$wsTrials.Range("C237").Value = 7.603698253631592
$wsTrials.Range("B238").Value = 4122
$wsTrials.Range("C238").Value = 2.264414548873901
$wsTrials.Range("B239").Value = 20437
$wsTrials.Range("C239").Value = 38.57143521308899
$wsTrials.Range("B240").Value = 8045
$wsTrials.Range("C240").Value = 6.867434978485107
$wsTrials.Range("B241").Value = 747
$wsTrials.Range("C241").Value = 0.3430116176605225
$wsTrials.Range("B242").Value = 1157
$wsTrials.Range("C242").Value = 0.3716356754302979
$wsTrials.Range("B243").Value = 6478
$wsTrials.Range("C243").Value = 4.722816944122314
$wsTrials.Range("B244").Value = 6809
$wsTrials.Range("C244").Value = 4.339688777923584
$wsTrials.Range("D244").Value = 1
$wsTrials.Range("B245").Value = 12019
$wsTrials.Range("C245").Value = 13.74567031860352
$wsTrials.Range("B246").Value = 1771
$wsTrials.Range("C246").Value = 0.9404408931732178
$wsTrials.Range("B247").Value = 7726
$wsTrials.Range("C247").Value = 7.326365232467651
$wsTrials.Range("B248").Value = 3887
$wsTrials.Range("C248").Value = 2.703219652175903
$wsTrials.Range("B249").Value = 11623
$wsTrials.Range("C249").Value = 13.80415940284729
$wsTrials.Range("B250").Value = 9558
$wsTrials.Range("C250").Value = 8.711540699005127
$wsTrials.Range("D250").Value = 1
$wsTrials.Range("B251").Value = 1406
$wsTrials.Range("C251").Value = 0.4971628189086914
$wsTrials.Range("B252").Value = 5346
$wsTrials.Range("C252").Value = 3.22569727897644
$wsTrials.Range("D252").Value = 1
$wsTrials.Range("B253").Value = 7336
$wsTrials.Range("C253").Value = 8.07194995880127
$wsTrials.Range("B254").Value = 929
$wsTrials.Range("C254").Value = 0.4008622169494629
$wsTrials.Range("B255").Value = 3108
$wsTrials.Range("C255").Value = 1.669046401977539
$wsTrials.Range("B256").Value = 987
$wsTrials.Range("C256").Value = 0.303565502166748
$wsTrials.Range("B257").Value = 4472
$wsTrials.Range("C257").Value = 3.398208856582642
$wsTrials.Range("B258").Value = 5488
$wsTrials.Range("C258").Value = 4.856003046035767
$wsTrials.Range("B259").Value = 1030
$wsTrials.Range("C259").Value = 0.4567477703094482
$wsTrials.Range("B260").Value = 4794
$wsTrials.Range("C260").Value = 2.370518684387207
$wsTrials.Range("B261").Value = 1360
$wsTrials.Range("C261").Value = 0.5222876071929932
$wsTrials.Range("C262").Value = 11.58659934997559
$wsTrials.Range("B263").Value = 5666
$wsTrials.Range("C263").Value = 4.741200923919678
$wsTrials.Range("B264").Value = 9316
$wsTrials.Range("C264").Value = 10.93622636795044
$wsTrials.Range("B265").Value = 5789
$wsTrials.Range("C265").Value = 4.467854022979736
$wsTrials.Range("B266").Value = 25053
$wsTrials.Range("C266").Value = 60.00302410125732
$wsTrials.Range("D266").Value = 0
$wsTrials.Range("B267").Value = 7993
$wsTrials.Range("C267").Value = 7.092535495758057
$wsTrials.Range("B268").Value = 4653
$wsTrials.Range("C268").Value = 3.052749872207642
$wsTrials.Range("D268").Value = 1
$wsTrials.Range("B269").Value = 2765
$wsTrials.Range("C269").Value = 1.416789293289185
$wsTrials.Range("D269").Value = 1
$wsTrials.Range("B270").Value = 2235
$wsTrials.Range("C270").Value = 0.8996939659118652
$wsTrials.Range("D270").Value = 1
$wsTrials.Range("B271").Value = 14262
$wsTrials.Range("C271").Value = 21.42769384384155
$wsTrials.Range("B272").Value = 1339
$wsTrials.Range("C272").Value = 0.4966344833374023
$wsTrials.Range("B273").Value = 15188
$wsTrials.Range("C273").Value = 23.69496846199036
$wsTrials.Range("B274").Value = 9504
$wsTrials.Range("C274").Value = 10.50045824050903
$wsTrials.Range("B275").Value = 3972
$wsTrials.Range("C275").Value = 2.626878976821899
$wsTrials.Range("B276").Value = 9588
$wsTrials.Range("C276").Value = 9.237613677978516
$wsTrials.Range("B277").Value = 6689
$wsTrials.Range("C277").Value = 6.771347761154175
$wsTrials.Range("B278").Value = 4112
$wsTrials.Range("C278").Value = 2.830191612243652
$wsTrials.Range("D278").Value = 1
$wsTrials.Range("B279").Value = 3248
$wsTrials.Range("C279").Value = 1.855069398880005
$wsTrials.Range("B280").Value = 3816
$wsTrials.Range("C280").Value = 2.355343341827393
$wsTrials.Range("B281").Value = 4850
$wsTrials.Range("C281").Value = 3.203450202941895
$wsTrials.Range("B282").Value = 2659
$wsTrials.Range("C282").Value = 1.279846668243408
$wsTrials.Range("B283").Value = 24527
$wsTrials.Range("C283").Value = 55.60902070999146
$wsTrials.Range("B284").Value = 4768
$wsTrials.Range("C284").Value = 3.208128690719604
$wsTrials.Range("B285").Value = 20440
$wsTrials.Range("C285").Value = 38.82545948028564
$wsTrials.Range("D285").Value = 1
$wsTrials.Range("B286").Value = 8797
$wsTrials.Range("C286").Value = 8.277260303497314
$wsTrials.Range("B287").Value = 2003
$wsTrials.Range("C287").Value = 0.8834660053253174
$wsTrials.Range("B288").Value = 8417
$wsTrials.Range("C288").Value = 7.813794851303101
$wsTrials.Range("B289").Value = 2907
$wsTrials.Range("C289").Value = 1.306690454483032
$wsTrials.Range("B290").Value = 26792
$wsTrials.Range("C290").Value = 60.00456547737122
$wsTrials.Range("D290").Value = 0
$wsTrials.Range("B291").Value = 7825
$wsTrials.Range("C291").Value = 6.064299821853638
$wsTrials.Range("B292").Value = 2062
$wsTrials.Range("C292").Value = 0.7566604614257812
$wsTrials.Range("D292").Value = 1
$wsTrials.Range("B293").Value = 1062
$wsTrials.Range("C293").Value = 0.3848831653594971
$wsTrials.Range("B294").Value = 3184
$wsTrials.Range("C294").Value = 1.201053380966187
$wsTrials.Range("B295").Value = 2417
$wsTrials.Range("C295").Value = 0.8760776519775391
$wsTrials.Range("B296").Value = 6241
$wsTrials.Range("C296").Value = 3.969612598419189
$wsTrials.Range("B297").Value = 4784
$wsTrials.Range("C297").Value = 2.708643198013306
$wsTrials.Range("B298").Value = 1820
$wsTrials.Range("C298").Value = 0.6968178749084473
$wsTrials.Range("B299").Value = 2265
$wsTrials.Range("C299").Value = 0.8960909843444824
$wsTrials.Range("B300").Value = 3153
$wsTrials.Range("C300").Value = 1.296282768249512
$wsTrials.Range("B301").Value = 1529
$wsTrials.Range("C301").Value = 0.5607259273529053
$wsSummary.Range("C2").Value = 0.3266666666666667
